$d = $word.ActiveDocument

# 1. "Regelhoogte: 1.428 (20px)" -> "Regelhoogte: 1.428" (drop the pixel hint)
$d.Content.Find.Execute("Regelhoogte: 1.428 (20px)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Regelhoogte: 1.428", 2) | Out-Null

# 2. Font update: "Lato" -> "Open Sans" (the actual edit behind the commit message)
$d.Content.Find.Execute("Lettertype: Lato (regular 400)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Lettertype: Open Sans (regular 400)", 2) | Out-Null
